$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 18: fill in a new time-log entry for "Arsi"
$ws.Range("A18").Value = "Arsi"
$ws.Range("B18").Value = "2024-02-12"
$ws.Range("C18").Formula = "=10+21/60"
$ws.Range("D18").Formula = "=10+42/60"
$ws.Range("F18").Value = "Button input working -> no even link"

# Restore the view: scroll back to top and move the selection to F19
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("F19").Select()
